$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue ($ws.Range("D2")) '309.60'
Set-TextValue ($ws.Range("E2")) '-0.31%'
Set-TextValue ($ws.Range("D3")) '37.21'
Set-TextValue ($ws.Range("E3")) '-0.91%'
Set-TextValue ($ws.Range("D4")) '5.124'
Set-TextValue ($ws.Range("E4")) '0.23%'
Set-TextValue ($ws.Range("D5")) '0.07848'
Set-TextValue ($ws.Range("E5")) '0.55%'
Set-TextValue ($ws.Range("B6")) 'KuCoinToken'
Set-TextValue ($ws.Range("C6")) 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue ($ws.Range("D6")) '8.271'
Set-TextValue ($ws.Range("E6")) '0.58%'
Set-TextValue ($ws.Range("B7")) 'FTXToken'
Set-TextValue ($ws.Range("C7")) 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue ($ws.Range("D7")) '1.882'
Set-TextValue ($ws.Range("E7")) '-1.36%'
Set-TextValue ($ws.Range("B8")) 'BTSEToken'
Set-TextValue ($ws.Range("C8")) 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue ($ws.Range("D8")) '2.999'
Set-TextValue ($ws.Range("E8")) '2.91%'
Set-TextValue ($ws.Range("B9")) 'MXToken'
Set-TextValue ($ws.Range("C9")) 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue ($ws.Range("D9")) '0.9231'
Set-TextValue ($ws.Range("E9")) '-0.34%'
Set-TextValue ($ws.Range("B10")) 'LiechtensteinCryptoassetsExchange'
Set-TextValue ($ws.Range("C10")) 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue ($ws.Range("D10")) '0.1100'
Set-TextValue ($ws.Range("E10")) '-7.93%'
Set-TextValue ($ws.Range("B11")) 'WazirX'
Set-TextValue ($ws.Range("C11")) 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue ($ws.Range("D11")) '0.1904'
Set-TextValue ($ws.Range("E11")) '0.09%'
Set-TextValue ($ws.Range("B12")) 'MandalaExchangeToken'
Set-TextValue ($ws.Range("C12")) 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue ($ws.Range("D12")) '0.08920'
Set-TextValue ($ws.Range("E12")) '-5.29%'
Set-TextValue ($ws.Range("B13")) 'BitrueCoin'
Set-TextValue ($ws.Range("C13")) 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue ($ws.Range("D13")) '0.03320'
Set-TextValue ($ws.Range("E13")) '-3.45%'
Set-TextValue ($ws.Range("B14")) 'BitMartToken'
Set-TextValue ($ws.Range("C14")) 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue ($ws.Range("D14")) '0.09606'
Set-TextValue ($ws.Range("E14")) '-0.13%'
Set-TextValue ($ws.Range("B15")) 'BitForexToken'
Set-TextValue ($ws.Range("C15")) 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue ($ws.Range("D15")) '0.001378'
Set-TextValue ($ws.Range("E15")) '1.03%'
Set-TextValue ($ws.Range("B16")) 'TigerCash'
Set-TextValue ($ws.Range("C16")) 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue ($ws.Range("D16")) '0.005995'
Set-TextValue ($ws.Range("E16")) '1.58%'
Set-TextValue ($ws.Range("B17")) 'LEO'
Set-TextValue ($ws.Range("C17")) 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue ($ws.Range("D17")) '3.392'
Set-TextValue ($ws.Range("E17")) '-4.00%'
Set-TextValue ($ws.Range("B18")) 'GateToken'
Set-TextValue ($ws.Range("C18")) 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue ($ws.Range("D18")) '4.392'
Set-TextValue ($ws.Range("E18")) '-0.20%'
Set-TextValue ($ws.Range("D19")) '0.3435'
Set-TextValue ($ws.Range("E19")) '0.31%'
Set-TextValue ($ws.Range("D20")) '6.362'
Set-TextValue ($ws.Range("E20")) '20.92%'
Set-TextValue ($ws.Range("D21")) '0.1298'
Set-TextValue ($ws.Range("E21")) '2.35%'
Set-TextValue ($ws.Range("D22")) '0.2409'
Set-TextValue ($ws.Range("E22")) '-7.01%'
Set-TextValue ($ws.Range("D23")) '0.04354'
Set-TextValue ($ws.Range("E23")) '-0.02%'
Set-TextValue ($ws.Range("D24")) '0.001199'
Set-TextValue ($ws.Range("E24")) '0.03%'
Set-TextValue ($ws.Range("D25")) '0.004286'
Set-TextValue ($ws.Range("E25")) '0.68%'
Set-TextValue ($ws.Range("D26")) '0.0001400'
Set-TextValue ($ws.Range("E26")) '7.60%'
Set-TextValue ($ws.Range("D27")) '0.0002899'
Set-TextValue ($ws.Range("D39")) '0.02170'
Set-TextValue ($ws.Range("E39")) '4.49%'
Set-TextValue ($ws.Range("D40")) '0.05024'
Set-TextValue ($ws.Range("E40")) '-1.09%'
Set-TextValue ($ws.Range("D41")) '0.007580'
Set-TextValue ($ws.Range("E41")) '-0.81%'
Set-TextValue ($ws.Range("D42")) '0.1357'
Set-TextValue ($ws.Range("E42")) '0.68%'
Set-TextValue ($ws.Range("D43")) '0.008512'
Set-TextValue ($ws.Range("E43")) '-6.77%'
Set-TextValue ($ws.Range("D44")) '0.002007'
Set-TextValue ($ws.Range("E44")) '-3.09%'
Set-TextValue ($ws.Range("D45")) '0.008100'
Set-TextValue ($ws.Range("E45")) '-5.91%'
Set-TextValue ($ws.Range("D46")) '0.00006536'
Set-TextValue ($ws.Range("E46")) '-2.63%'
Set-TextValue ($ws.Range("D47")) '0.00000000750'
Set-TextValue ($ws.Range("E47")) '-0.07%'
Set-TextValue ($ws.Range("D48")) '0.003295'
Set-TextValue ($ws.Range("E48")) '13.05%'
Set-TextValue ($ws.Range("D49")) '0.001443'
Set-TextValue ($ws.Range("E49")) '20.16%'
Set-TextValue ($ws.Range("D50")) '0.00002099'
Set-TextValue ($ws.Range("E50")) '-0.07%'
Set-TextValue ($ws.Range("D51")) '0.0001999'
Set-TextValue ($ws.Range("E51")) '-0.07%'
